$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$wsTypes = $wb.Worksheets.Item("Data types")

# --- "Sheet1": add the new test_junos device row ---
$ws1.Range("A10").Value = "192.168.20.244"
$ws1.Range("C10").Value = "home-qfx"
$ws1.Range("B10").Value = "test_junos"
$ws1.Range("D10").Value = "qfx"
$ws1.Range("E10").Value = "backbone"

# --- "Sheet1": a few existing rows had their group reset to "none" ---
$ws1.Range("B2").Value = "none"
$ws1.Range("B7").Value = "none"
$ws1.Range("B8").Value = "none"

# --- "Data types" sheet: insert the new "test_junos" group before "none" ---
# Current groups list (E6:E8): junos, cisco, none
# New groups list (E6:E9):     junos, cisco, test_junos, none
$wsTypes.Range("E8").Value = "test_junos"
$wsTypes.Range("E9").Value = "none"

# Extend the "groups" named range so it covers the new row too.
$wb.Names.Item("groups").RefersTo = "='Data types'!`$E`$6:`$E`$9"
